$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1748.6781
$ws.Range("I15").Value = 1748.6781
$ws.Range("K15").Value = 5246.0343
$ws.Range("M15").Value = -5077.0343

$ws.Range("H116").Value = 3000
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H135").Value = 852.69446
$ws.Range("I135").Value = 373.44116
$ws.Range("J135").Value = 9000
$ws.Range("K135").Value = 3360.97044
$ws.Range("L135").Value = 81000
$ws.Range("M135").Value = -825.9704400000001
$ws.Range("N135").Value = -86070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 539204.6
$ws.Range("I32").Value = 642214.1
$ws.Range("K32").Value = 642214.1
$ws.Range("M32").Value = -641927.1

$ws.Range("H45").Value = 3566.3333
$ws.Range("I45").Value = 3349.5
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 3349.5
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -2972.5
$ws.Range("N45").Value = -4754

$ws.Range("H74").Value = 1113.2858
$ws.Range("I74").Value = 766.2692
$ws.Range("J74").Value = 2115.7778
$ws.Range("K74").Value = 766.2692
$ws.Range("L74").Value = 2115.7778
$ws.Range("M74").Value = 107.7308
$ws.Range("N74").Value = -3863.7778

$ws.Range("H77").Value = 1113.2858
$ws.Range("I77").Value = 766.2692
$ws.Range("J77").Value = 2115.7778
$ws.Range("K77").Value = 3831.346
$ws.Range("L77").Value = 10578.889
$ws.Range("M77").Value = 536.6540000000005
$ws.Range("N77").Value = -19314.889

$ws.Range("H122").Value = 201680
$ws.Range("I122").Value = 250850
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 752550
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -750100
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 3268.3704
$ws.Range("I132").Value = 2089.5898
$ws.Range("J132").Value = 6333.2
$ws.Range("K132").Value = 6268.769400000001
$ws.Range("L132").Value = 18999.6
$ws.Range("M132").Value = -3738.769400000001
$ws.Range("N132").Value = -24059.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 72205.664
$ws.Range("J9").Value = 72205.664
$ws.Range("L9").Value = 72205.664
$ws.Range("N9").Value = -72541.664

$ws.Range("H58").Value = 2127.9412
$ws.Range("I58").Value = 2058.8572
$ws.Range("J58").Value = 2176.3
$ws.Range("K58").Value = 2058.8572
$ws.Range("L58").Value = 2176.3
$ws.Range("M58").Value = -1855.8572
$ws.Range("N58").Value = -2582.3

$ws.Range("H92").Value = 50533
$ws.Range("J92").Value = 50533
$ws.Range("L92").Value = 50533
$ws.Range("N92").Value = -55525

$ws.Range("H97").Value = 39800
$ws.Range("J97").Value = 39800
$ws.Range("L97").Value = 39800
$ws.Range("N97").Value = -41782

$ws.Range("H105").Value = 2999
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2999
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2999
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -6493

$ws.Range("H132").Value = 9261816
$ws.Range("I132").Value = 2502.375
$ws.Range("K132").Value = 7507.125
$ws.Range("M132").Value = -4977.125

$ws.Range("H134").Value = 5267.7334
$ws.Range("I134").Value = 5293.0415
$ws.Range("J134").Value = 5166.5
$ws.Range("K134").Value = 15879.1245
$ws.Range("L134").Value = 15499.5
$ws.Range("M134").Value = -13344.1245
$ws.Range("N134").Value = -20569.5

$ws.Range("H136").Value = 2127.9412
$ws.Range("I136").Value = 2058.8572
$ws.Range("J136").Value = 2176.3
$ws.Range("K136").Value = 6176.571599999999
$ws.Range("L136").Value = 6528.900000000001
$ws.Range("M136").Value = -3626.571599999999
$ws.Range("N136").Value = -11628.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1571.9333
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1571.9333
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 4715.7999
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -5291.7999

$ws.Range("H80").Value = 6794
$ws.Range("J80").Value = 6794
$ws.Range("L80").Value = 20382
$ws.Range("N80").Value = -22254

$ws.Range("H83").Value = 6794
$ws.Range("J83").Value = 6794
$ws.Range("L83").Value = 61146
$ws.Range("N83").Value = -70506

$ws.Range("H102").Value = 3444.4443
$ws.Range("J102").Value = 3444.4443
$ws.Range("L102").Value = 10333.3329
$ws.Range("N102").Value = -15201.3329

$ws.Range("H108").Value = 1000
$ws.Range("I108").Value = 1000
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3000
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -120
$ws.Range("N108").ClearContents()

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H117").Value = 518.5714
$ws.Range("J117").Value = 866.6667
$ws.Range("L117").Value = 2600.0001
$ws.Range("N117").Value = -9484.000100000001

$ws.Range("H122").Value = 4528.346
$ws.Range("I122").Value = 377.3125
$ws.Range("J122").Value = 11170
$ws.Range("K122").Value = 3395.8125
$ws.Range("L122").Value = 100530
$ws.Range("M122").Value = -945.8125
$ws.Range("N122").Value = -105430

$ws.Range("H125").Value = 2013.375
$ws.Range("I125").Value = 2020
$ws.Range("J125").Value = 2012.9333
$ws.Range("K125").Value = 6060
$ws.Range("L125").Value = 6038.7999
$ws.Range("M125").Value = -1140
$ws.Range("N125").Value = -15878.7999

$ws.Range("H133").Value = 14709
$ws.Range("I133").Value = 9540
$ws.Range("J133").Value = 19016.5
$ws.Range("K133").Value = 28620
$ws.Range("L133").Value = 57049.5
$ws.Range("M133").Value = -23560
$ws.Range("N133").Value = -67169.5

$ws.Range("H137").Value = 8342636
$ws.Range("I137").Value = 33353920
$ws.Range("K137").Value = 100061760
$ws.Range("M137").Value = -100056660

$ws.Range("H138").Value = 2409.8125
$ws.Range("I138").Value = 825.5714
$ws.Range("J138").Value = 13499.5
$ws.Range("K138").Value = 2476.7142
$ws.Range("L138").Value = 40498.5
$ws.Range("M138").Value = 2663.2858
$ws.Range("N138").Value = -50778.5

$ws.Range("H139").Value = 4479.4
$ws.Range("I139").Value = 1799.2307
$ws.Range("J139").Value = 9456.857
$ws.Range("K139").Value = 5397.6921
$ws.Range("L139").Value = 28370.571
$ws.Range("M139").Value = -257.6921000000002
$ws.Range("N139").Value = -38650.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4274522
$ws.Range("I136").Value = 907.3226
$ws.Range("J136").Value = 20834780
$ws.Range("K136").Value = 2721.9678
$ws.Range("L136").Value = 62504340
$ws.Range("M136").Value = -171.9677999999999
$ws.Range("N136").Value = -62509440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H132").Value = 4506764.5
$ws.Range("I132").Value = 2501.7727
$ws.Range("J132").Value = 11113017
$ws.Range("K132").Value = 7505.3181
$ws.Range("L132").Value = 33339051
$ws.Range("M132").Value = -4975.3181
$ws.Range("N132").Value = -33344111

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 2139
$ws.Range("I136").Value = 1576.1471
$ws.Range("J136").Value = 3733.75
$ws.Range("K136").Value = 4728.4413
$ws.Range("L136").Value = 11201.25
$ws.Range("M136").Value = -2178.4413
$ws.Range("N136").Value = -16301.25
